# Update TPM-derived LR-pair metrics (Anxa2-Robo4) with recomputed values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (ECs -> ECs)
$ws.Range("G2").Value = 32.09557633333333
$ws.Range("H2").Value = 96.28672900000001
$ws.Range("I2").Value = 0.1656600924295661
$ws.Range("J2").Value = 0.1656600924295661
$ws.Range("M2").Value = 39.327127
$ws.Range("N2").Value = 117.981381
$ws.Range("O2").Value = 0.9923865713449503
$ws.Range("P2").Value = 0.9923865713449502
$ws.Range("Q2").Value = 1262.226806599194
$ws.Range("R2").Value = 11360.04125939275
$ws.Range("S2").Value = 0.1643988511348646
$ws.Range("T2").Value = 0.1643988511348646

# Row 3 (ECs -> FAPs)
$ws.Range("G3").Value = 32.09557633333333
$ws.Range("H3").Value = 96.28672900000001
$ws.Range("I3").Value = 0.1656600924295661
$ws.Range("J3").Value = 0.1656600924295661
$ws.Range("O3").Value = 0.001455135597170125
$ws.Range("P3").Value = 0.001455135597170125
$ws.Range("Q3").Value = 1.850802107787111
$ws.Range("R3").Value = 16.657218970084
$ws.Range("S3").Value = 0.0002410578975247547
$ws.Range("T3").Value = 0.0002410578975247547

# Row 4 (ECs -> MuSCs)
$ws.Range("G4").Value = 32.09557633333333
$ws.Range("H4").Value = 96.28672900000001
$ws.Range("I4").Value = 0.1656600924295661
$ws.Range("J4").Value = 0.1656600924295661
$ws.Range("O4").Value = 0.00615829305787961
$ws.Range("P4").Value = 0.006158293057879609
$ws.Range("Q4").Value = 7.832797021844666
$ws.Range("R4").Value = 70.495173196602
$ws.Range("S4").Value = 0.001020183397176691
$ws.Range("T4").Value = 0.001020183397176691

# Row 5 (FAPs -> ECs)
$ws.Range("I5").Value = 0.6938590312037638
$ws.Range("J5").Value = 0.6938590312037638
$ws.Range("M5").Value = 39.327127
$ws.Range("N5").Value = 117.981381
$ws.Range("O5").Value = 0.9923865713449503
$ws.Range("P5").Value = 0.9923865713449502
$ws.Range("Q5").Value = 5286.77399814144
$ws.Range("R5").Value = 47580.96598327296
$ws.Range("S5").Value = 0.6885763849730321
$ws.Range("T5").Value = 0.688576384973032

# Row 6 (FAPs -> FAPs)
$ws.Range("I6").Value = 0.6938590312037638
$ws.Range("J6").Value = 0.6938590312037638
$ws.Range("O6").Value = 0.001455135597170125
$ws.Range("P6").Value = 0.001455135597170125
$ws.Range("S6").Value = 0.001009658975722573
$ws.Range("T6").Value = 0.001009658975722573

# Row 7 (FAPs -> MuSCs)
$ws.Range("I7").Value = 0.6938590312037638
$ws.Range("J7").Value = 0.6938590312037638
$ws.Range("O7").Value = 0.00615829305787961
$ws.Range("P7").Value = 0.006158293057879609
$ws.Range("S7").Value = 0.00427298725500921
$ws.Range("T7").Value = 0.00427298725500921

# Row 8 (MuSCs -> ECs)
$ws.Range("I8").Value = 0.1404808763666701
$ws.Range("J8").Value = 0.1404808763666701
$ws.Range("M8").Value = 39.327127
$ws.Range("N8").Value = 117.981381
$ws.Range("O8").Value = 0.9923865713449503
$ws.Range("P8").Value = 0.9923865713449502
$ws.Range("Q8").Value = 1070.376850356697
$ws.Range("R8").Value = 9633.391653210276
$ws.Range("S8").Value = 0.1394113352370537
$ws.Range("T8").Value = 0.1394113352370536

# Row 9 (MuSCs -> FAPs)
$ws.Range("I9").Value = 0.1404808763666701
$ws.Range("J9").Value = 0.1404808763666701
$ws.Range("O9").Value = 0.001455135597170125
$ws.Range("P9").Value = 0.001455135597170125
$ws.Range("S9").Value = 0.0002044187239227971
$ws.Range("T9").Value = 0.0002044187239227971

# Row 10 (MuSCs -> MuSCs)
$ws.Range("I10").Value = 0.1404808763666701
$ws.Range("J10").Value = 0.1404808763666701
$ws.Range("O10").Value = 0.00615829305787961
$ws.Range("P10").Value = 0.006158293057879609
$ws.Range("S10").Value = 0.0008651224056937084
$ws.Range("T10").Value = 0.0008651224056937083
